# LSTM_mv.xlsx daily-update edit script
# Appends newly observed actual values / predictions for several currency &
# commodity pairs, matching the workflow of the original workbook (each
# weekday a new prediction row is completed and a new blank placeholder row
# is appended for the following day's "actual" value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "D1_USD" (sheet1) : rows 77-82
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("D1_USD")

$ws1.Range("A77").Value = 45265
$ws1.Range("B77").Value = 3.9918

$ws1.Range("A78").Value = 45266
$ws1.Range("B78").Value = 4.0019920000000004
$ws1.Range("C78").Value = 3.9952830000000001

$ws1.Range("A79").Value = 45267
$ws1.Range("B79").Value = 4.0191990000000004
$ws1.Range("C79").Value = 4.0261930000000001

$ws1.Range("A80").Value = 45268
$ws1.Range("B80").Value = 4.0191990000000004
$ws1.Range("C80").Value = 4.0016084000000003

$ws1.Range("A81").Value = 45271
$ws1.Range("B81").Value = 4.01816
$ws1.Range("C81").Value = 4.0067434000000004

$ws1.Range("C82").Value = 4.0336319999999999

$ws1.Range("D69:D81").Formula = "=B69-C69"
$ws1.Range("E69:E81").Formula = "=IF(D69<0,1,0)"

# restore date-format look of the new date cells (copy format from A76)
$ws1.Range("A76").Copy()
$ws1.Range("A77:A81").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("C83").Select()

# ---------------------------------------------------------------------
# Sheet "D1_EUR" (sheet3) : rows 351-356
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("D1_EUR")

$ws3.Range("A351").Value = 45265
$ws3.Range("B351").Value = 4.32761

$ws3.Range("A352").Value = 45266
$ws3.Range("B352").Value = 4.3194900000000001
$ws3.Range("C352").Value = 4.3752639999999996

$ws3.Range("A353").Value = 45267
$ws3.Range("B353").Value = 4.3272300000000001
$ws3.Range("C353").Value = 4.3687243000000002

$ws3.Range("A354").Value = 45268
$ws3.Range("B354").Value = 4.3272300000000001
$ws3.Range("C354").Value = 4.3779199999999996

$ws3.Range("A355").Value = 45271
$ws3.Range("B355").Value = 4.3254599999999996
$ws3.Range("C355").Value = 4.3646349999999998

$ws3.Range("C356").Value = 4.3341146000000004

$ws3.Range("D346:D355").Formula = "=B346-C346"
$ws3.Range("E346:E355").Formula = "=IF(D346<0,1,0)"

$ws3.Range("A350").Copy()
$ws3.Range("A351:A355").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("C357").Select()

# ---------------------------------------------------------------------
# Sheet "D5_EUR" (sheet5) : rows 44-46, B turns from "Nan" text to numbers
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("D5_EUR")

$ws5.Range("B44").Value = 4.32761
$ws5.Range("B45").Value = 4.3194900000000001
$ws5.Range("B46").Value = 4.3272300000000001

$ws5.Range("B47").Select()

# ---------------------------------------------------------------------
# Sheet "D1_OIL" (sheet6) : backfill history (rows 2-13), complete
# existing rows (14-22) and append new rows (23-28)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("D1_OIL")

# existing rows 2-10 move down to rows 14-22 -- read off the old values
# first (still addressable at their old location before we overwrite them)
$oldA = @()
$oldB = @()
$oldC = @()
for ($i = 2; $i -le 10; $i++) {
    $oldA += $ws6.Cells.Item($i, 1).Value2
    $oldB += $ws6.Cells.Item($i, 2).Value2
    $oldC += $ws6.Cells.Item($i, 3).Value2
}

$newHistory = @(
    @(45233, 82.459998999999996, 79.287199999999999),
    @(45236, 80.819999999999993, 78.5989),
    @(45237, 77.370002999999997, 79.054299999999998),
    @(45238, 75.330001999999993, 77.506699999999995),
    @(45239, 75.739998, 76.052700000000002),
    @(45240, 77.169998000000007, 75.321899999999999),
    @(45244, 78.260002, 78.169200000000004),
    @(45245, 78.260002, 76.024100000000004),
    @(45246, 76.660004000000001, 74.524199999999993),
    @(45247, 72.900002000000001, 75.210099999999997),
    @(45250, 77.599997999999999, 73.762),
    @(45251, 77.769997000000004, 75.093100000000007)
)

$tail = @(
    @(73.526899999999998),
    @(75.0822),
    @(71.261700000000005)
)

$newTail = @(
    @(45265, 72.319999999999993, 73.881699999999995),
    @(45266, 69.379997000000003, 71.1126),
    @(45267, 69.339995999999999, 73.859700000000004),
    @(45268, 71.23, 71.077100000000002),
    @(45271, 71.319999999999993, 68.771699999999996)
)

# write the back-filled historical rows into rows 2-13
$r = 2
foreach ($row in $newHistory) {
    $ws6.Cells.Item($r, 1).Value = $row[0]
    $ws6.Cells.Item($r, 2).Value = $row[1]
    $ws6.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# re-write the old rows (2-10 previously) into their new location (14-22),
# filling in column C (the "Day+1 Prediction" placeholder that had just
# been entered the following day) that had not been available before
$missingC = @(72.002200000000002, 70.452699999999993, 76.697699999999998)
for ($i = 0; $i -le 8; $i++) {
    $rowNum = 14 + $i
    $ws6.Cells.Item($rowNum, 1).Value = $oldA[$i]
    $ws6.Cells.Item($rowNum, 2).Value = $oldB[$i]
    if ($oldC[$i]) {
        $ws6.Cells.Item($rowNum, 3).Value = $oldC[$i]
    }
}
$ws6.Cells.Item(14, 3).Value = 72.002200000000002
$ws6.Cells.Item(15, 3).Value = 70.452699999999993
$ws6.Cells.Item(16, 3).Value = 76.697699999999998

# append the new tail rows (23-27) and the trailing placeholder row (28)
$r = 23
foreach ($row in $newTail) {
    $ws6.Cells.Item($r, 1).Value = $row[0]
    $ws6.Cells.Item($r, 2).Value = $row[1]
    $ws6.Cells.Item($r, 3).Value = $row[2]
    $r++
}
$ws6.Cells.Item(28, 3).Value = 69.378200000000007

# fill in the Difference / Ratio formulas for every data row (2-27)
$ws6.Range("D2:D27").Formula = "=B2-C2"
$ws6.Range("E2:E27").Formula = "=D2/C2"

$ws6.Range("E1").Formula = "=AVERAGE(D2:D301)"

# restore formats: column A date format, column B / E numeric formats,
# and the yellow highlight on E1
$ws6.Range("A2").Copy()
$ws6.Range("A3:A27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Range("B2").Copy()
$ws6.Range("B3:B27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Range("E2").Copy()
$ws6.Range("E3:E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("C225").Copy()
$ws6.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Activate()
$ws6.Range("C29").Select()
